$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 197.57143
$ws.Range("I2").Value = 197.57143
$ws.Range("K2").Value = 197.57143
$ws.Range("M2").Value = -84.57142999999999

$ws.Range("H5").Value = 192.5
$ws.Range("I5").Value = 85
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 85
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = -530

$ws.Range("H9").Value = 244.83333
$ws.Range("I9").Value = 244.83333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 244.83333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -75.83332999999999
$ws.Range("N9").ClearContents()

$ws.Range("H28").Value = 819
$ws.Range("I28").Value = 795.625
$ws.Range("K28").Value = 795.625
$ws.Range("M28").Value = -310.625

$ws.Range("H31").Value = 459.6
$ws.Range("I31").Value = 459.6
$ws.Range("K31").Value = 1378.8
$ws.Range("M31").Value = -1148.8

$ws.Range("H43").Value = 2823.077
$ws.Range("I43").Value = 1450
$ws.Range("J43").Value = 3433.3333
$ws.Range("K43").Value = 1450
$ws.Range("L43").Value = 3433.3333
$ws.Range("M43").Value = -1381
$ws.Range("N43").Value = -3571.3333

$ws.Range("H70").Value = 7935.4116
$ws.Range("J70").Value = 11500.182
$ws.Range("L70").Value = 34500.546
$ws.Range("N70").Value = -35040.546

$ws.Range("H73").Value = 7935.4116
$ws.Range("J73").Value = 11500.182
$ws.Range("L73").Value = 34500.546
$ws.Range("N73").Value = -36372.546

$ws.Range("H98").Value = 1424.909
$ws.Range("I98").Value = 1424.909
$ws.Range("K98").Value = 1424.909
$ws.Range("M98").Value = 73.09099999999989

$ws.Range("H122").Value = 1424.909
$ws.Range("I122").Value = 1424.909
$ws.Range("K122").Value = 4274.727000000001
$ws.Range("M122").Value = -1824.727000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 109

$ws.Range("H61").Value = 2283.4
$ws.Range("I61").Value = 2283.4
$ws.Range("K61").Value = 2283.4
$ws.Range("M61").Value = -2071.4

$ws.Range("H97").Value = 2781.5557
$ws.Range("I97").Value = 2154.25
$ws.Range("K97").Value = 2154.25
$ws.Range("M97").Value = -1658.25

$ws.Range("H132").Value = 1592.5714
$ws.Range("I132").Value = 1408
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 4224
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -1694
$ws.Range("N132").Value = -13160

$ws.Range("H136").Value = 2283.4
$ws.Range("I136").Value = 2283.4
$ws.Range("K136").Value = 6850.200000000001
$ws.Range("M136").Value = -4300.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 112

$ws.Range("H22").Value = 3196.375
$ws.Range("I22").Value = 3594.8333
$ws.Range("J22").Value = 2001
$ws.Range("K22").Value = 3594.8333
$ws.Range("L22").Value = 2001
$ws.Range("M22").Value = -3421.8333
$ws.Range("N22").Value = -2347

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H134").Value = 7054.7144
$ws.Range("I134").Value = 5781.143
$ws.Range("K134").Value = 17343.429
$ws.Range("M134").Value = -14808.429

$ws.Range("H135").Value = 55999
$ws.Range("J135").Value = 55999
$ws.Range("L135").Value = 55999
$ws.Range("N135").Value = -66139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 245.08333
$ws.Range("I7").Value = 175.85715
$ws.Range("J7").Value = 342
$ws.Range("K7").Value = 175.85715
$ws.Range("L7").Value = 342
$ws.Range("M7").Value = -62.85714999999999
$ws.Range("N7").Value = -568

$ws.Range("H31").Value = 7253
$ws.Range("I31").Value = 6341.8335
$ws.Range("K31").Value = 6341.8335
$ws.Range("M31").Value = -6046.8335

$ws.Range("H34").Value = 7253
$ws.Range("I34").Value = 6341.8335
$ws.Range("K34").Value = 6341.8335
$ws.Range("M34").Value = -6139.8335

$ws.Range("H50").Value = 29500
$ws.Range("J50").Value = 29333.334
$ws.Range("L50").Value = 29333.334
$ws.Range("N50").Value = -30583.334

$ws.Range("H51").Value = 24333.334
$ws.Range("J51").Value = 24333.334
$ws.Range("L51").Value = 24333.334
$ws.Range("N51").Value = -25805.334

$ws.Range("H59").Value = 28333.334
$ws.Range("I59").Value = 15000
$ws.Range("K59").Value = 15000
$ws.Range("M59").Value = -13855

$ws.Range("H60").Value = 22400
$ws.Range("J60").Value = 25500
$ws.Range("L60").Value = 25500
$ws.Range("N60").Value = -26522

$ws.Range("H61").Value = 24333.334
$ws.Range("J61").Value = 24333.334
$ws.Range("L61").Value = 24333.334
$ws.Range("N61").Value = -25029.334

$ws.Range("H68").Value = 42499.5
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51497

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

$ws.Range("H71").Value = 42499.5
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -157485

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H86").Value = 4099.75
$ws.Range("I86").Value = 4099.75
$ws.Range("K86").Value = 4099.75
$ws.Range("M86").Value = -2976.75

$ws.Range("H89").Value = 4099.75
$ws.Range("I89").Value = 4099.75
$ws.Range("K89").Value = 20498.75
$ws.Range("M89").Value = -14882.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 656655.5
$ws.Range("I4").Value = 678984.9399999999
$ws.Range("J4").Value = 500349.5
$ws.Range("K4").Value = 2036954.82
$ws.Range("L4").Value = 1501048.5
$ws.Range("M4").Value = -2036842.82
$ws.Range("N4").Value = -1501272.5

$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -18630

$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -20184

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 132.77777
$ws.Range("I2").Value = 148.75
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 148.75
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = -35.75
$ws.Range("N2").Value = -231

$ws.Range("H97").Value = 1745.6666
$ws.Range("J97").Value = 691.6667
$ws.Range("L97").Value = 691.6667
$ws.Range("N97").Value = -1683.6667

$ws.Range("H126").Value = 5331.5
$ws.Range("I126").Value = 5497.25
$ws.Range("K126").Value = 16491.75
$ws.Range("M126").Value = -14021.75

$ws.Range("H132").Value = 3338.2
$ws.Range("I132").Value = 3319.5
$ws.Range("K132").Value = 9958.5
$ws.Range("M132").Value = -7428.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2813.7856
$ws.Range("J22").Value = 2889
$ws.Range("L22").Value = 2889
$ws.Range("N22").Value = -3479

$ws.Range("H27").Value = 2813.7856
$ws.Range("J27").Value = 2889
$ws.Range("L27").Value = 2889
$ws.Range("N27").Value = -3103

$ws.Range("H46").Value = 1628.5714
$ws.Range("I46").Value = 1733
$ws.Range("K46").Value = 1733
$ws.Range("M46").Value = -1545

$ws.Range("H55").Value = 1721.8889
$ws.Range("I55").Value = 1440
$ws.Range("K55").Value = 1440
$ws.Range("M55").Value = -1267

$ws.Range("H93").Value = 4062.75
$ws.Range("I93").Value = 4143.143
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 4143.143
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = -2895.143
$ws.Range("N93").Value = -5996

$ws.Range("H122").Value = 6600
$ws.Range("I122").Value = 8200
$ws.Range("K122").Value = 24600
$ws.Range("M122").Value = -22150

